$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Status column (K) for rows 8, 16, 22 from "Failed" to "passed" ---
$ws.Range("K8").Value = "passed"
$ws.Range("K16").Value = "passed"
$ws.Range("K22").Value = "passed"

# --- Fill in the new "Bug ID" column (L) values.
# Order matters: it controls the order new shared strings are appended,
# so write them in the same sequence the original author did.
$ws.Range("L5").Value = "Car_SearchBug_02"
$ws.Range("L2").Value = "Car_SearchBug_01"
$ws.Range("L10").Value = "Car_SearchBug_03"
$ws.Range("L9").Value = "Car_SearchBug_04"
$ws.Range("L12").Value = "Car_SearchBug_05"
$ws.Range("L19").Value = "Car_SearchBug_06"
$ws.Range("L25").Value = "Car_SearchBug_07"

# Apply wrap-text-only formatting (no centering) to the new Bug ID cells,
# matching the new style used by the author (wrapText, default alignment).
# Format the first cell directly, then copy/paste-special its format onto
# the rest so only a single new style entry is created (matches the diff's
# cellXfs count going from 5 to 6, instead of accumulating stray styles).
$firstBugCell = $ws.Range("L2")
$firstBugCell.HorizontalAlignment = 1
$firstBugCell.VerticalAlignment = -4107
$firstBugCell.WrapText = $true

$firstBugCell.Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("L12").PasteSpecial(-4122)
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("L25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: scroll so row 15 / column G is the top-left visible cell,
# and select L15 (best effort - some view state may not round-trip). ---
$win = $excel.ActiveWindow
$excel.Goto($ws.Range("G15"), $true)
$win.ScrollRow = 15
$win.ScrollColumn = 7
$ws.Range("L15").Select()
